# Update "Plan de test 2048" worksheet with actual test execution results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update two scenario/expected-result texts whose wording changed ---
$ws.Range("A16").Value = "Boutons radios 2048 et 4x4 cochés par défaut"
$ws.Range("B16").Value = "Les boutons radio doivent être cochés au chargement de la page menu."

$ws.Range("A22").Value = "Grille de jeu avec deux blocs par défaut"
$ws.Range("B22").Value = "Afficher correctement une grille de jeu avec deux blocs au chargement de la page de jeu."

# --- Row 16: Conforme ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = $null
$ws.Range("I16").Value = 1

# --- Row 17: Conforme ---
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = $null
$ws.Range("I17").Value = 2

# --- Row 18: Conforme ---
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = $null
$ws.Range("I18").Value = 1

# --- Row 19: Conforme ---
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = $null
$ws.Range("I19").Value = 1

# --- Row 20: Conforme ---
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = $null
$ws.Range("I20").Value = 1

# --- Row 22: Conforme ---
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = $null
$ws.Range("I22").Value = 1

# --- Row 23: Conforme ---
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = $null
$ws.Range("I23").Value = 1

# --- Row 24: Conforme (nombre d'essai left blank, as in source) ---
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = $null

# --- Row 25: Conforme ---
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = $null
$ws.Range("I25").Value = 1

# --- Row 26: Non-Conforme - "Pas d’animations" / "Développer les animations" ---
$ws.Range("E26").Value = $null
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = "Pas d’animations"
$ws.Range("H26").Value = "Développer les animations"
$ws.Range("I26").Value = 2

# --- Row 27: Conforme ---
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = $null
$ws.Range("I27").Value = 1

# --- Row 28: Conforme ---
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = $null
$ws.Range("I28").Value = 1

# --- Row 29: Conforme ---
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = $null
$ws.Range("I29").Value = 1

# --- Row 30: Non-Conforme - "Victoire pas détectée" ---
$ws.Range("E30").Value = $null
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = "Victoire pas détectée"
$ws.Range("H30").Value = "Arrêter la partie quand il y a un 2048 sur la grille pour le mode classique"
$ws.Range("I30").Value = 2

# --- Row 31: Non-Conforme - "fichier csv non trouvé" ---
$ws.Range("E31").Value = $null
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = "fichier csv non trouvé"
$ws.Range("I31").Value = 1

# --- Row 32: Conforme ---
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = $null
$ws.Range("I32").Value = 1

# --- Row 33: Non-Conforme - "Bouton non présent" / "Ajouter le bouton" ---
$ws.Range("E33").Value = $null
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = "Bouton non présent"
$ws.Range("H33").Value = "Ajouter le bouton"
$ws.Range("I33").Value = 1

# --- Row 34: Non-Conforme - Responsive not adapted ---
$ws.Range("E34").Value = $null
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = "Responsive pas adapté à toutes les tailles d’écran"
$ws.Range("H34").Value = "Rapetisser la grille et les blocs en fonction de l’écran"
$ws.Range("I34").Value = 1

# --- Row 36: Conforme ---
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = $null
$ws.Range("I36").Value = 1

# --- Row 37: Conforme ---
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = $null
$ws.Range("I37").Value = 1

# --- Row 38: Non-Conforme - "Pas d’animations" / "Développer les animations" ---
$ws.Range("E38").Value = $null
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = "Pas d’animations"
$ws.Range("H38").Value = "Développer les animations"
$ws.Range("I38").Value = 1

# --- Row 39: Conforme ---
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = $null
$ws.Range("I39").Value = 1

# --- Row 40: Non-Conforme - "Pas d’animations" / "Développer les animations" ---
$ws.Range("E40").Value = $null
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = "Pas d’animations"
$ws.Range("H40").Value = "Développer les animations"
$ws.Range("I40").Value = 1

# --- Row 41: Conforme ---
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = $null
$ws.Range("I41").Value = 1

# --- Row 42: Non-Conforme - "Pas d’animations" / "Développer les animations" ---
$ws.Range("E42").Value = $null
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = "Pas d’animations"
$ws.Range("H42").Value = "Développer les animations"
$ws.Range("I42").Value = 1

# --- Row 43: Conforme ---
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = $null
$ws.Range("I43").Value = 1

# --- Row 44: Non-Conforme - "Bouton non présent" / "Ajouter le bouton" ---
$ws.Range("E44").Value = $null
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = "Bouton non présent"
$ws.Range("H44").Value = "Ajouter le bouton"
$ws.Range("I44").Value = 1

# --- Row 45: Non-Conforme - Responsive not adapted (table version) ---
$ws.Range("E45").Value = $null
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = "Responsive pas adapté à toutes les tailles d’écran"
$ws.Range("H45").Value = "Rapetisser le tableau en fonction de l’écran"
$ws.Range("I45").Value = 1

# --- Restore cursor/selection roughly where the author left it ---
$ws.Range("H36").Select()
